# Daily Status Tracker update
# - Renumber the "SlNo" column (A) so it is a contiguous sequence again
#   (rows 5-23 previously skipped 4, 9 and 24).
# - Move the active selection / view to B6 (also resets the frozen-pane
#   top-left cell back to row 2, matching a fresh view after scrolling
#   back to the top of the data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracker")

$slNo = @{
    5  = 4
    6  = 5
    7  = 6
    8  = 7
    9  = 8
    10 = 9
    11 = 10
    12 = 11
    13 = 12
    14 = 13
    15 = 14
    16 = 15
    17 = 16
    18 = 17
    19 = 18
    20 = 19
    21 = 20
    22 = 21
    23 = 22
}

foreach ($row in $slNo.Keys) {
    $ws.Cells.Item($row, 1).Value = $slNo[$row]
}

# Restore selection to B6 (updates activeCell/sqref and resets the pane's
# topLeftCell back to the top of the scrollable area).
$ws.Range("B6").Select()
